# Fix: populate the "Emission factor - Description" and "Emission factor - Source"
# columns (V and W) on the eCRF_3 sheet with their corresponding template
# placeholders, for every data row that already has the GHG columns (Y:AB)
# populated (rows 2-65). Rows 66-71 are "blank" template rows (no GHG values
# filled in either) and are intentionally left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eCRF_3")

for ($r = 2; $r -le 65; $r++) {
    $vCell = $ws.Range("V$r")
    $vCell.Value2 = "{{emission_factor_description}}"
    $vCell.NumberFormat = "@"

    $wCell = $ws.Range("W$r")
    $wCell.Value2 = "{{emission_factor_source}}"
    $wCell.NumberFormat = "@"
}
